$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.479.37"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.897.75"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'527.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").Value = "'142.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.16%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "2.898.33"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").Value = "'5.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("D12").Value = "'0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "3.404.99"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "60.490.60"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "'22.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").Value = "2.906.00"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("D19").Value = "'5.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").Value = "'11.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'363.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'64.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "3.023.84"
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("D26").Value = "'0.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'7.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.77%  "
$ws.Range("D30").Value = "0.0₃0859"
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("D32").Value = "'1.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Value = "'19.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").Value = "'145.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.40%  "
$ws.Range("D35").Value = "'4.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("D36").Value = "'5.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.98%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("D38").Value = "'1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.16%  "
$ws.Range("D39").Value = "'37.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("D40").Value = "'1.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("D41").Value = "2.323.52"
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").Value = "'3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.31%  "
$ws.Range("D43").Value = "'0.643"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "'20.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.18%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'4.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "'0.0234"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("D49").Value = "'0.0933"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("D50").Value = "'10.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'249.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.36%  "
